# Generate Report for Handoff
#
# This CI run re-generated the localization status report: the batch of
# files that were previously "Ready for handoff" / "Handback transform
# failed" just had a fresh handoff performed, so their "Latest Handoff
# Date(time)" timestamps on the Overview sheet and on each locale sheet
# (zh-cn, de-de) are refreshed to the new handoff run's timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet - column D = "Latest Handoff Date"
# Rows 7, 10-16 correspond to files that were handed off in this run.
$newHandoffDateOverview = "2016-20-11 22:20:17"
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $overview.Range("D$r").Value = $newHandoffDateOverview
}

# zh-cn sheet - column E = "Latest Handoff Datetime"
$newHandoffDateTimeZhCn = "2016-03-11 22:20:13"
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $zhcn.Range("E$r").Value = $newHandoffDateTimeZhCn
}

# de-de sheet - column E = "Latest Handoff Datetime"
$newHandoffDateTimeDeDe = "2016-03-11 22:20:17"
foreach ($r in @(7, 10, 11, 12, 13, 14, 15, 16)) {
    $dede.Range("E$r").Value = $newHandoffDateTimeDeDe
}
